$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates extracted from the target diff.
# Each entry only lists the columns (B=Coin, C=Link, D=Price, E=Volume) that change for that row.
$updates = @(
    @{ row=2; D="29.083.31" }
    @{ row=3; D="1.830.06"; E="  -0.35%  " }
    @{ row=4; D="0.9963"; E="  -0.63%  " }
    @{ row=5; D="242.97"; E="  -0.04%  " }
    @{ row=6; D="0.6288"; E="  +0.38%  " }
    @{ row=7; D="0.9985"; E="  -0.31%  " }
    @{ row=8; D="0.07527"; E="  -0.57%  " }
    @{ row=9; D="0.2925"; E="  -0.09%  " }
    @{ row=10; D="23.30"; E="  +3.37%  " }
    @{ row=11; D="0.07667"; E="  -0.87%  " }
    @{ row=12; D="1.837.13"; E="  +0.34%  " }
    @{ row=13; D="5.012"; E="  +1.05%  " }
    @{ row=14; D="0.6692"; E="  +0.86%  " }
    @{ row=15; D="82.86"; E="  +0.16%  " }
    @{ row=16; D="0.000009403"; E="  -6.37%  " }
    @{ row=17; D="5.989"; E="  -0.98%  " }
    @{ row=18; D="29.081.51"; E="  +0.12%  " }
    @{ row=19; D="2.080.52"; E="  +0.15%  " }
    @{ row=20; D="12.57"; E="  +1.76%  " }
    @{ row=21; D="223.28"; E="  -1.52%  " }
    @{ row=22; D="0.9999"; E="  -0.28%  " }
    @{ row=23; D="7.117"; E="  -1.09%  " }
    @{ row=24; D="0.9972"; E="  -0.57%  " }
    @{ row=25; D="159.65"; E="  +0.54%  " }
    @{ row=26; D="0.1395"; E="  +1.12%  " }
    @{ row=27; D="8.498"; E="  -0.06%  " }
    @{ row=28; D="17.88"; E="  -0.23%  " }
    @{ row=29; D="1.493"; E="  -0.34%  " }
    @{ row=30; E="  +10.23%  " }
    @{ row=31; D="4.152"; E="  +1.19%  " }
    @{ row=32; D="4.099"; E="  +2.13%  " }
    @{ row=33; D="1.206"; E="  +0.92%  " }
    @{ row=34; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7410"; E="  +1.06%  " }
    @{ row=35; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.834"; E="  -0.33%  " }
    @{ row=36; D="1.138"; E="  +0.01%  " }
    @{ row=37; D="2.663"; E="  -1.10%  " }
    @{ row=38; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.224.92"; E="  -1.33%  " }
    @{ row=39; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.758"; E="  -0.25%  " }
    @{ row=40; D="0.01779"; E="  -0.18%  " }
    @{ row=41; D="6.498"; E="  +2.56%  " }
    @{ row=42; D="0.8876"; E="  -1.13%  " }
    @{ row=43; D="0.9988"; E="  -0.22%  " }
    @{ row=45; B="BabyDogeCoin"; C="https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D="0.00000000126"; E="  +1.48%  " }
    @{ row=46; B="RocketPoolETH"; C="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D="1.977.49"; E="  +0.12%  " }
    @{ row=47; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="65.79"; E="  +2.44%  " }
    @{ row=48; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="0.5082"; E="  -0.77%  " }
    @{ row=49; B="XinFinNetwork"; C="https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"; D="0.07591"; E="  +12.90%  " }
    @{ row=50; D="0.4067"; E="  +0.69%  " }
    @{ row=51; D="8.981"; E="  +1.38%  " }
)

# A cell whose new text looks like a plain number (e.g. '0.9985') would be
# auto-coerced to a Number by Excel's normal typed-input behaviour. The
# source data keeps these as text (prices like '1.830.06' use '.' as a
# thousands separator), so for any numeric-looking value we first mark the
# cell as Text (NumberFormat '@') before writing it, preserving the string.
$numericLike = '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$'

foreach ($u in $updates) {
    $r = $u.row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) {
        $dVal = $u.D
        if ($dVal -match $numericLike) {
            $ws.Cells.Item($r, 4).NumberFormat = "@"
        }
        $ws.Cells.Item($r, 4).Value = $dVal
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}